$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1301
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 2502
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 7506
$ws.Range("M6").Value = -188
$ws.Range("N6").Value = -7730
$ws.Range("H121").Value = 1499.6666
$ws.Range("J121").Value = 1499.6666
$ws.Range("L121").Value = 4498.9998
$ws.Range("N121").Value = -7992.9998
$ws.Range("H137").Value = 1482.6451
$ws.Range("I137").Value = 1112.4667
$ws.Range("J137").Value = 1829.6875
$ws.Range("K137").Value = 3337.4001
$ws.Range("L137").Value = 5489.0625
$ws.Range("M137").Value = -787.4000999999998
$ws.Range("N137").Value = -10589.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 745.36365
$ws.Range("I2").Value = 488.42856
$ws.Range("J2").Value = 1195
$ws.Range("K2").Value = 488.42856
$ws.Range("L2").Value = 1195
$ws.Range("M2").Value = -375.42856
$ws.Range("N2").Value = -1421
$ws.Range("H32").Value = 3238.682
$ws.Range("I32").Value = 2027.3572
$ws.Range("K32").Value = 2027.3572
$ws.Range("M32").Value = -1740.3572
$ws.Range("H63").Value = 10000
$ws.Range("I63").Value = 10000
$ws.Range("K63").Value = 10000
$ws.Range("M63").Value = -9314
$ws.Range("H66").Value = 10000
$ws.Range("I66").Value = 10000
$ws.Range("K66").Value = 50000
$ws.Range("M66").Value = -46568
$ws.Range("H116").Value = 745.36365
$ws.Range("I116").Value = 488.42856
$ws.Range("J116").Value = 1195
$ws.Range("K116").Value = 488.42856
$ws.Range("L116").Value = 1195
$ws.Range("M116").Value = 1805.57144
$ws.Range("N116").Value = -5783
$ws.Range("H132").Value = 2503.5557
$ws.Range("I132").Value = 1790.5714
$ws.Range("K132").Value = 5371.7142
$ws.Range("M132").Value = -2841.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 745.36365
$ws.Range("I3").Value = 488.42856
$ws.Range("J3").Value = 1195
$ws.Range("K3").Value = 488.42856
$ws.Range("L3").Value = 1195
$ws.Range("M3").Value = -374.42856
$ws.Range("N3").Value = -1423
$ws.Range("H107").Value = 2024.8334
$ws.Range("I107").Value = 2024.8334
$ws.Range("K107").Value = 2024.8334
$ws.Range("M107").Value = -104.8334
$ws.Range("H108").Value = 64985
$ws.Range("J108").Value = 64985
$ws.Range("L108").Value = 64985
$ws.Range("N108").Value = -72665
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 522
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 3490
$ws.Range("J31").Value = 6203.375
$ws.Range("L31").Value = 6203.375
$ws.Range("N31").Value = -6793.375
$ws.Range("H34").Value = 3490
$ws.Range("J34").Value = 6203.375
$ws.Range("L34").Value = 6203.375
$ws.Range("N34").Value = -6607.375
$ws.Range("H113").Value = 522
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 549.25
$ws.Range("J7").Value = 1000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224
$ws.Range("H12").Value = 119.5
$ws.Range("J12").Value = 135.625
$ws.Range("L12").Value = 406.875
$ws.Range("N12").Value = -752.875
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H23").Value = 128.2
$ws.Range("J23").Value = 181
$ws.Range("L23").Value = 543
$ws.Range("N23").Value = -1013
$ws.Range("H37").Value = 99999
$ws.Range("J37").Value = 99999
$ws.Range("L37").Value = 299997
$ws.Range("N37").Value = -300221
$ws.Range("H92").Value = 340
$ws.Range("J92").Value = 340
$ws.Range("L92").Value = 1020
$ws.Range("N92").Value = -3516
$ws.Range("H122").Value = 675.3333
$ws.Range("I122").Value = 298.85715
$ws.Range("K122").Value = 2689.71435
$ws.Range("M122").Value = -239.7143499999997
$ws.Range("H131").Value = 792.0700000000001
$ws.Range("J131").Value = 795.0404
$ws.Range("L131").Value = 2385.1212
$ws.Range("N131").Value = -12465.1212
$ws.Range("H136").Value = 500000000
$ws.Range("I136").Value = 500000000
$ws.Range("K136").Value = 1500000000
$ws.Range("M136").Value = -1499994900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6711043
$ws.Range("J11").Value = 3109232
$ws.Range("L11").Value = 3109232
$ws.Range("N11").Value = -3109510
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H70").Value = 5313.25
$ws.Range("I70").Value = 5584.3335
$ws.Range("K70").Value = 5584.3335
$ws.Range("M70").Value = -5314.3335
$ws.Range("H73").Value = 5313.25
$ws.Range("I73").Value = 5584.3335
$ws.Range("K73").Value = 5584.3335
$ws.Range("M73").Value = -4648.3335
$ws.Range("H113").Value = 1522.2
$ws.Range("I113").Value = 1403.6666
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 1403.6666
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = 766.3334
$ws.Range("N113").Value = -6040
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H126").Value = 69605.13
$ws.Range("I126").Value = 3062.3
$ws.Range("J126").Value = 202690.8
$ws.Range("K126").Value = 9186.900000000001
$ws.Range("L126").Value = 608072.3999999999
$ws.Range("M126").Value = -6716.900000000001
$ws.Range("N126").Value = -613012.3999999999
$ws.Range("H132").Value = 5126.385
$ws.Range("I132").Value = 4181.2104
$ws.Range("K132").Value = 12543.6312
$ws.Range("M132").Value = -10013.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5553.2666
$ws.Range("I7").Value = 3057.4285
$ws.Range("J7").Value = 7737.125
$ws.Range("K7").Value = 3057.4285
$ws.Range("L7").Value = 7737.125
$ws.Range("M7").Value = -2945.4285
$ws.Range("N7").Value = -7961.125
$ws.Range("H126").Value = 5553.2666
$ws.Range("I126").Value = 3057.4285
$ws.Range("J126").Value = 7737.125
$ws.Range("K126").Value = 9172.2855
$ws.Range("L126").Value = 23211.375
$ws.Range("M126").Value = -6702.2855
$ws.Range("N126").Value = -28151.375
$ws.Range("H132").Value = 2307.5557
$ws.Range("I132").Value = 1949.5
$ws.Range("J132").Value = 2409.8572
$ws.Range("K132").Value = 5848.5
$ws.Range("L132").Value = 7229.571599999999
$ws.Range("M132").Value = -3318.5
$ws.Range("N132").Value = -12289.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 18000
$ws.Range("J20").Value = 18000
$ws.Range("L20").Value = 18000
$ws.Range("N20").Value = -18480
$ws.Range("H122").Value = 129211.7
$ws.Range("I122").Value = 129211.7
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 387635.1
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -385185.1
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3265.25
$ws.Range("I132").Value = 3044.5
$ws.Range("K132").Value = 9133.5
$ws.Range("M132").Value = -6603.5
